# Prevent Word's "smart quotes" autocorrect from turning straight
# quotes into curly quotes when we set text programmatically.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Edit 1: insert a new "Meta description" paragraph right after the
# title heading ("Play Bohemian Bazaar Slot for Free - Review 2021").
# -------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Find out what we like and don' + [char]0x27 + 't like about Bohemian Bazaar, a High 5 Games slot with a gypsy theme and a bonus game. Play for free today!</w:t></w:r>' +
           '</w:p>'
$metaPara.Range.InsertXML($metaXml)

# -------------------------------------------------------------------
# Edit 2: at the end of the document, drop the duplicated bold title
# paragraph and rewrite the italic paragraph with the new image-prompt
# text (keeping its existing run/paragraph formatting).
# -------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($count - 1)
$boldTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)
$italicRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End)
$italicRange.Text = "Please create an image that features a happy Maya warrior with glasses in cartoon style, fitting the theme of the game ""Bohemian Bazaar"". The Maya warrior should be wearing traditional clothing with a headpiece, possibly adorned with feathers or flowers, and be depicted in a joyful pose with a big smile on their face. The background of the image should include colorful tents, carriages, and banners that reflect the bohemian and gypsy culture of the game. The overall style should be cheerful and vibrant, with bold colors and fun details that capture the essence of the game."

Write-Output "Edit complete."
Write-Output "Paragraph 1: $($d.Paragraphs.Item(1).Range.Text)"
Write-Output "Paragraph 2: $($d.Paragraphs.Item(2).Range.Text)"
Write-Output "Last paragraph: $($d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)"
